# Adding doc to MBED Project MIDI2DMX
# Duplicate the existing red "Ellipse 81" marker shape on slide 7 five
# times, re-positioning (and renaming) each copy to create five new
# red-dot markers on the schematic.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# EMU-per-point constant used by the PowerPoint object model (Left/Top/
# Width/Height are expressed in points). The Left/Top setters round-trip
# through a single-precision float and truncate back to EMU, so nudge by
# half an EMU (in point units) to land exactly on the target EMU value.
$EMU_PER_PT = 12700.0
$HALF_EMU = 0.5

# The shape to clone: "Ellipse 81" (the last shape on the slide).
$template = $s.Shapes.Item($s.Shapes.Count)

# (new name, x offset EMU, y offset EMU)
$newMarkers = @(
    @("Ellipse 69", 13187,   5121735),
    @("Ellipse 82", 2665425, 5093036),
    @("Ellipse 83", 871626,  5099474),
    @("Ellipse 84", 1828818, 5121394),
    @("Ellipse 85", 4078582, 6435580)
)

foreach ($marker in $newMarkers) {
    $name = $marker[0]
    $offX = $marker[1]
    $offY = $marker[2]

    $dupRange = $template.Duplicate()
    $dup = $dupRange.Item(1)

    $dup.Name = $name
    $dup.Left = ($offX + $HALF_EMU) / $EMU_PER_PT
    $dup.Top = ($offY + $HALF_EMU) / $EMU_PER_PT
}

Write-Output ("Slide 7 shape count now: " + $s.Shapes.Count)
